$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A (so old A,B,C shift to C,D,E)
$ws.Range("A1:B1").EntireColumn.Insert()

# Row 1 headers
$ws.Range("A1").Value = "Level"
$ws.Range("C1").Value = "Subsection(agriculture product)"
$ws.Range("B1").Value = "State"

# Row 2 values
$ws.Range("A2").Value = "Kabupaten"
$ws.Range("B2").Value = "Aceh"

# New columns F and G for Start Year / End Year values
$ws.Range("F2").Value = 1970
$ws.Range("G2").Value = 2024

# Style A2 with special font
$ws.Range("A2").Font.Name = "Courier New"
$ws.Range("A2").Font.Size = 7
$ws.Range("A2").Font.Color = 14935011

# Column widths (values chosen so the runtime's internal pixel rounding
# reproduces the exact target XML "width" attributes)
$ws.Columns.Item(1).ColumnWidth = 7.571428571428571
$ws.Columns.Item(2).ColumnWidth = 5.142857142857143
$ws.Columns.Item(3).ColumnWidth = 29.285714285714285
$ws.Columns.Item(4).ColumnWidth = 8.142857142857142
$ws.Columns.Item(5).ColumnWidth = 9.857142857142858
$ws.Columns.Item(6).ColumnWidth = 9.142857142857142
$ws.Columns.Item(7).ColumnWidth = 8.428571428571429

# Selection
$ws.Range("B4").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
